$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting the existing rows 102-110 down to 103-111.
$ws.Rows.Item(102).Insert()

# Populate the new row 102 with the new weekly record (same market/category template
# as the surrounding rows, with its own date/volume/price figures).
$ws.Cells.Item(102, 1).Value = 5
$ws.Cells.Item(102, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(102, 3).Value = "Maule"
$ws.Cells.Item(102, 4).Value = 44461
$ws.Cells.Item(102, 5).Value = 7
$ws.Cells.Item(102, 6).Value = 100112017
$ws.Cells.Item(102, 7).Value = "Apio"
$ws.Cells.Item(102, 8).Value = "Americana (o)"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 500
$ws.Cells.Item(102, 11).Value = 9000
$ws.Cells.Item(102, 12).Value = 9000
$ws.Cells.Item(102, 13).Value = 9000
$ws.Cells.Item(102, 14).Value = "$/docena de matas"
$ws.Cells.Item(102, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(102, 16).Value = 1500
$ws.Cells.Item(102, 17).Value = 6
$ws.Cells.Item(102, 18).Value = "Hortaliza"
